# Add description about java installation on Mac OS X (Lion)
# Slide 14, Shape 2 ("TextShape 2"), 5th paragraph:
#   "Software : Java 1.6.x" -> "Software : Java 1.6.x " (trailing space added)
#   followed by new orange-colored runs describing the java install steps.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$para = $tr.Paragraphs(5)

# Orange accent color used elsewhere on the slide (ff420e), expressed as a
# COM RGB long (0x00BBGGRR = R + G*256 + B*65536).
$orange = 934655

# --- Update the existing run's text in place: add a trailing space ---
$start = $para.Start
$len = $para.Text.Length
$firstRun = $tr.Characters($start, $len)
$firstRun.Text = $firstRun.Text + " "

# --- Append the new orange runs describing the Java install steps ---
$newTexts = @(
    "Mac OS X (Lion)",
    "をお使いの方は",
    "java",
    "のインストールが必要です。ターミナルを起動し、",
    "java",
    "と入力することで、",
    "java",
    "をインストールしてください。"
)

foreach ($t in $newTexts) {
    $beforeStart = $para.Start
    $beforeLen = $para.Text.Length
    $null = $para.InsertAfter($t)
    $newStart = $beforeStart + $beforeLen
    $newRun = $tr.Characters($newStart, $t.Length)
    $newRun.Font.Color.RGB = $orange
}

Write-Host "Final paragraph text:" $para.Text
